# Incremental Extract Id Log.xlsx - apply "Upload before pulling files from BC Cert DB Server" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New text cells -> these introduce new shared strings. Order matters: the
#    diff's sharedStrings.xml gains "Supplier" (index 8) before
#    "Same next Id" (index 9), so write B11 before G9.
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Supplier"
$ws.Range("G9").Value = "Same next Id"

# ---------------------------------------------------------------------------
# 2. Fill in numeric values added to existing rows (columns C/E)
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = 2217696
$ws.Range("C8").Value = 1120842
$ws.Range("C9").Value = 2321664
$ws.Range("C10").Value = 1161211

# ---------------------------------------------------------------------------
# 3. Rows 13-16: repeat of the "header" block seen in rows 8-11, reusing the
#    same date style (s="1") as column A in that block and the same shared
#    strings for column B.
# ---------------------------------------------------------------------------
$ws.Range("A8").Copy()
$ws.Range("A13:A15").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A13").Value = 43596.813194444447
$ws.Range("B13").Value = "Item"

$ws.Range("A14").Value = 43596.809027777781
$ws.Range("B14").Value = "Supplier Item"

$ws.Range("A15").Value = 43596.813194444447
$ws.Range("B15").Value = "Special"

$ws.Range("B16").Value = "Supplier"

# ---------------------------------------------------------------------------
# 4. Rows 18-20 and 22-24: new date-only blocks using a short-date style
#    (built-in numFmtId 14) instead of the existing date+time custom format.
# ---------------------------------------------------------------------------
$ws.Range("A18").NumberFormat = "mm-dd-yy"
$ws.Range("A18").Value = 43635
$ws.Range("B18").Value = "Item"
$ws.Range("C18").Value = 1124379

$ws.Range("A18").Copy()
$ws.Range("A19:A20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A22:A24").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A19").Value = 43635
$ws.Range("B19").Value = "Supplier Item"
$ws.Range("C19").Value = 2396774

$ws.Range("A20").Value = 43635
$ws.Range("B20").Value = "Special"
$ws.Range("C20").Value = 1161916

$ws.Range("A22").Value = 43659
$ws.Range("B22").Value = "Item"
$ws.Range("C22").Value = 1125083

$ws.Range("A23").Value = 43659
$ws.Range("B23").Value = "Supplier Item"
$ws.Range("C23").Value = 2401993

$ws.Range("A24").Value = 43659
$ws.Range("B24").Value = "Special"
$ws.Range("C24").Value = 1163221

# ---------------------------------------------------------------------------
# 5. Reposition / resize the screenshot picture: it moves from hanging below
#    the data (near F14) up to the top-right of the sheet (near E1).
# ---------------------------------------------------------------------------
function Get-ColX($ws, $colIndex0, $offEmu) {
    $total = 0.0
    for ($i = 1; $i -le $colIndex0; $i++) {
        $total += $ws.Columns.Item($i).Width
    }
    $total += $offEmu / 12700.0
    return $total
}

function Get-RowY($ws, $rowIndex0, $offEmu) {
    $total = 0.0
    for ($i = 1; $i -le $rowIndex0; $i++) {
        $total += $ws.Rows.Item($i).Height
    }
    $total += $offEmu / 12700.0
    return $total
}

$fromX = Get-ColX $ws 4 746760
$fromY = Get-RowY $ws 0 22860
$toX = Get-ColX $ws 21 441960
$toY = Get-RowY $ws 8 167640

$shp = $ws.Shapes.Item("Picture 1")
$shp.Left = $fromX
$shp.Top = $fromY
$shp.Width = $toX - $fromX
$shp.Height = $toY - $fromY

# ---------------------------------------------------------------------------
# 6. Update the view state: scrolled so row 6 is at the top, with A25 as the
#    active selection (cursor left below the newly entered data).
# ---------------------------------------------------------------------------
$excel.Windows.Item(1).ScrollRow = 6
$excel.Windows.Item(1).ScrollColumn = 1
$ws.Range("A25").Select()

Write-Output "edit applied"
